# 8.2.1.xlsx update: extend the "Labor productivity" table with the 2022
# data column (S) and refresh the previously-estimated 2020/2021 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New year column: S4 = 2022 -------------------------------------------
# Clone R4's formatting (right-aligned "year header" style with borders) so
# the new header cell matches its neighbours instead of getting a brand new
# style slot.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S4").Value = 2022

# --- Revise the 2020 figure (Q5) -------------------------------------------
# Previously an unformatted "General" estimate (90.6); now a finalized
# figure carrying the "0.0" number format used by the other recent years.
$ws.Range("Q5").Value = 91.892815141492093
$ws.Range("Q5").NumberFormat = "0.0"

# --- Revise the 2021 figure (R5) --------------------------------------------
$ws.Range("R5").Value = 101.53074848578628

# --- New 2022 figure (S5) ---------------------------------------------------
# Clone R5's formatting (the "0.0" number format + borders) for the new cell.
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S5").Value = 109.27053140096621

# --- Move the active selection ----------------------------------------------
[void]$ws.Range("T5").Select()
